# Update column G ("K") values on Sheet1 for rows 2-22.
# The source data regenerated "K" (strikeouts) from a corrected computation
# (K instead of Strike#) and these are the newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 4
    3  = 6
    4  = 9
    5  = 8
    6  = 7
    7  = 8
    8  = 3
    9  = 9
    10 = 6
    11 = 1
    12 = 5
    13 = 3
    14 = 9
    15 = 4
    16 = 3
    17 = 3
    18 = 5
    19 = 6
    20 = 3
    21 = 2
    22 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
